$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 207, shifting existing rows 207:329 down to 208:330
$ws.Rows(207).Insert()

# Populate the newly inserted row 207 with the new data record
$ws.Cells.Item(207, 1).Value = 3
$ws.Cells.Item(207, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(207, 3).Value = "Coquimbo"
$ws.Cells.Item(207, 4).Value = 44777
$ws.Cells.Item(207, 5).Value = 5
$ws.Cells.Item(207, 6).Value = 100112001
$ws.Cells.Item(207, 7).Value = "Berenjena"
$ws.Cells.Item(207, 8).Value = "Sin especificar"
$ws.Cells.Item(207, 9).Value = "Primera"
$ws.Cells.Item(207, 10).Value = 125
$ws.Cells.Item(207, 11).Value = 8500
$ws.Cells.Item(207, 12).Value = 9000
$ws.Cells.Item(207, 13).Value = 8740
$ws.Cells.Item(207, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(207, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(207, 16).Value = 146
$ws.Cells.Item(207, 17).Value = 60
$ws.Cells.Item(207, 18).Value = "Hortaliza"

# Match date cell style used by the rest of the D column (style index 2 in original file)
$ws.Cells.Item(207, 4).Style = $ws.Cells.Item(208, 4).Style
